$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ============================================================
# Sheet 1: "Indicadores"
# ============================================================

# First, copy the label-column style (bold, centered, bordered) from an
# existing labeled row (A9) down onto the new rows (15-23) that are being
# appended at the bottom of the sheet, so they match the existing look.
$ws1.Range("A9").Copy()
$ws1.Range("A15:A23").PasteSpecial(-4122)

# Row 2: Endividamento Total
$ws1.Range("A2").Value = "Endividamento Total"
$ws1.Range("B2").Value = 3.61068036787363
$ws1.Range("C2").Value = 3.65978964066897
$ws1.Range("D2").Value = 3.590901228701167
$ws1.Range("E2").Value = 3.547669704768785
$ws1.Range("F2").Value = 3.535321483732147

# Row 3: Dependência Financeira (%)
$ws1.Range("A3").Value = "Dependência Financeira (%)"
$ws1.Range("B3").Value = 100
$ws1.Range("C3").Value = 100
$ws1.Range("D3").Value = 100
$ws1.Range("E3").Value = 100
$ws1.Range("F3").Value = 100

# Row 4: Independência Financeira (%)
$ws1.Range("A4").Value = "Independência Financeira (%)"
$ws1.Range("B4").Value = 27.69561130078404
$ws1.Range("C4").Value = 27.32397482324178
$ws1.Range("D4").Value = 27.84816223869519
$ws1.Range("E4").Value = 28.18751696799164
$ws1.Range("F4").Value = 28.28597072717488

# Row 5: Capital Terceiros LP
$ws1.Range("A5").Value = "Capital Terceiros LP"
$ws1.Range("B5").Value = 2.61068036787363
$ws1.Range("C5").Value = 2.65978964066897
$ws1.Range("D5").Value = 2.590901228701166
$ws1.Range("E5").Value = 2.547669704768785
$ws1.Range("F5").Value = 2.535321483732147

# Row 6: Imobilização Capital Permanente
$ws1.Range("A6").Value = "Imobilização Capital Permanente"
$ws1.Range("B6").Value = 0.5436158758912167
$ws1.Range("C6").Value = 0.5841007837308343
$ws1.Range("D6").Value = 0.6332860718162319
$ws1.Range("E6").Value = 0.6642000495078371
$ws1.Range("F6").Value = 0.7170454445870662

# Row 7: Imobilização PL
$ws1.Range("A7").Value = "Imobilização PL"
$ws1.Range("B7").Value = 1.509805436284917
$ws1.Range("C7").Value = 1.542184752907285
$ws1.Range("D7").Value = 1.552688270253863
$ws1.Range("E7").Value = 1.555304830765864
$ws1.Range("F7").Value = 1.588187828365613

# Row 8: Imobilização Recursos Não Correntes
$ws1.Range("A8").Value = "Imobilização Recursos Não Correntes"
$ws1.Range("B8").Value = 0.506399173702311
$ws1.Range("C8").Value = 0.5413007481956502
$ws1.Range("D8").Value = 0.5737451953735432
$ws1.Range("E8").Value = 0.592655594514736
$ws1.Range("F8").Value = 0.619949725220023

# Row 9: Giro do Imobilizado
$ws1.Range("A9").Value = "Giro do Imobilizado"
$ws1.Range("B9").Value = 1.707850277633619
$ws1.Range("C9").Value = 2.214638267166006
$ws1.Range("D9").Value = 2.517534106021249
$ws1.Range("E9").Value = 2.912558117819684
$ws1.Range("F9").Value = 2.796014152384233

# Row 10: Composição Endividamento (%)
$ws1.Range("A10").Value = "Composição Endividamento (%)"
$ws1.Range("B10").Value = 72.30438869921596
$ws1.Range("C10").Value = 72.67602517675823
$ws1.Range("D10").Value = 72.15183776130478
$ws1.Range("E10").Value = 71.81248303200837
$ws1.Range("F10").Value = 71.71402927282512

# Row 11: Liquidez Geral
$ws1.Range("A11").Value = "Liquidez Geral"
$ws1.Range("B11").Value = 0.8443156610764809
$ws1.Range("C11").Value = 0.8386409868815722
$ws1.Range("D11").Value = 0.8430252061048135
$ws1.Range("E11").Value = 0.8477923465624001
$ws1.Range("F11").Value = 0.8528273445922694

# Row 12: Composição Endividamento LP (%)
$ws1.Range("A12").Value = "Composição Endividamento LP (%)"
$ws1.Range("B12").Value = 49.22448312178091
$ws1.Range("C12").Value = 44.81874187039704
$ws1.Range("D12").Value = 40.42985109600502
$ws1.Range("E12").Value = 37.81696668132022
$ws1.Range("F12").Value = 34.36477862424982

# Row 13: Participação ANC no Ativo (%)
$ws1.Range("A13").Value = "Participação ANC no Ativo (%)"
$ws1.Range("B13").Value = 41.81498450315775
$ws1.Range("C13").Value = 42.13861736122597
$ws1.Range("D13").Value = 43.2395148561486
$ws1.Range("E13").Value = 43.84018130761216
$ws1.Range("F13").Value = 44.92343442240518

# Row 14: Margem Operacional (%)
$ws1.Range("A14").Value = "Margem Operacional (%)"
$ws1.Range("B14").Value = -8.558666115731766
$ws1.Range("C14").Value = 4.652544000903408
$ws1.Range("D14").Value = -2.369369379742759
$ws1.Range("E14").Value = 5.831401193626308
$ws1.Range("F14").Value = 10.64843742016398

# Row 15: Margem Líquida (%)
$ws1.Range("A15").Value = "Margem Líquida (%)"
$ws1.Range("B15").Value = -18.30911419157669
$ws1.Range("C15").Value = -1.180509406717026
$ws1.Range("D15").Value = -4.465961137152987
$ws1.Range("E15").Value = 3.003981116115012
$ws1.Range("F15").Value = 5.430842226554105

# Row 16: ROA (%)
$ws1.Range("A16").Value = "ROA (%)"
$ws1.Range("B16").Value = -2.803384448797868
$ws1.Range("C16").Value = 1.773398941590173
$ws1.Range("D16").Value = -1.008247929558513
$ws1.Range("E16").Value = 2.927808167028303
$ws1.Range("F16").Value = 5.151783021070909

# Row 17: ROE (%)
$ws1.Range("A17").Value = "ROE (%)"
$ws1.Range("B17").Value = -21.65374177608858
$ws1.Range("C17").Value = -1.646802552087117
$ws1.Range("D17").Value = -6.824219169318738
$ws1.Range("E17").Value = 5.350693497321513
$ws1.Range("F17").Value = 9.288974679429796

# Row 18: GAO (Alavancagem Operacional)
$ws1.Range("A18").Value = "GAO (Alavancagem Operacional)"
$ws1.Range("C18").Value = -16.67343635890709
$ws1.Range("D18").Value = 103.0702857715295
$ws1.Range("E18").Value = -33.23052961155219
$ws1.Range("F18").Value = 4.65579857716681

# Row 19: GAF (Alavancagem Financeira)
$ws1.Range("A19").Value = "GAF (Alavancagem Financeira)"
$ws1.Range("C19").Value = 0.5824369222021166
$ws1.Range("D19").Value = -1.81640609962399
$ws1.Range("E19").Value = 0.4677093428249402
$ws1.Range("F19").Value = 0.9827286115496234

# Row 20: GAT (Alavancagem Total)
$ws1.Range("A20").Value = "GAT (Alavancagem Total)"
$ws1.Range("C20").Value = -9.71122495541471
$ws1.Range("D20").Value = -187.217495765394
$ws1.Range("E20").Value = -15.54222916634379
$ws1.Range("F20").Value = 4.575386471393852

# Row 21: Var. % Receita
$ws1.Range("A21").Value = "Var. % Receita"
$ws1.Range("C21").Value = 9.569885048300474
$ws1.Range("D21").Value = -1.457105385130242
$ws1.Range("E21").Value = 11.24872542319428
$ws1.Range("F21").Value = 29.19169630181118

# Row 22: Var. % EBIT
$ws1.Range("A22").Value = "Var. % EBIT"
$ws1.Range("C22").Value = -159.5628693148944
$ws1.Range("D22").Value = -150.1842684446086
$ws1.Range("E22").Value = -373.8011032676776
$ws1.Range("F22").Value = 135.9106581070581

# Row 23: Var. % Lucro Líquido
$ws1.Range("A23").Value = "Var. % Lucro Líquido"
$ws1.Range("C23").Value = -92.93530650150568
$ws1.Range("D23").Value = 272.7956212703539
$ws1.Range("E23").Value = -174.8302683565631
$ws1.Range("F23").Value = 133.5632923363448

# ============================================================
# Sheet 2: "Dados Base"
# ============================================================

# Row 2: AtivoCirculante
$ws2.Range("A2").Value = "AtivoCirculante"
$ws2.Range("B2").Value = 40549746.68045199
$ws2.Range("C2").Value = 37967991.53685346
$ws2.Range("D2").Value = 32876200.42487912
$ws2.Range("E2").Value = 30670519.513
$ws2.Range("F2").Value = 40326906

# Row 3: AtivoNaoCirculante
$ws2.Range("A3").Value = "AtivoNaoCirculante"
$ws2.Range("B3").Value = 29141300.63507662
$ws2.Range("C3").Value = 27650888.97605425
$ws2.Range("D3").Value = 25044728.79473388
$ws2.Range("E3").Value = 23942405.2206
$ws2.Range("F3").Value = 32892812

# Row 4: AtivoImobilizado
$ws2.Range("A4").Value = "AtivoImobilizado"
$ws2.Range("B4").Value = 13366071.25367216
$ws2.Range("C4").Value = 11293848.89331861
$ws2.Range("D4").Value = 9790271.2601236
$ws2.Range("E4").Value = 9414354.1516
$ws2.Range("F4").Value = 12669526

# Row 5: AtivoRealizavelLP
$ws2.Range("A5").Value = "AtivoRealizavelLP"
$ws2.Range("B5").Value = 1995054.146192828
$ws2.Range("C5").Value = 2026121.285431994
$ws2.Range("D5").Value = 2354678.50797864
$ws2.Range("E5").Value = 2578961.465
$ws2.Range("F5").Value = 4454043

# Row 6: AtivoTotal
$ws2.Range("A6").Value = "AtivoTotal"
$ws2.Range("B6").Value = 69691047.31552862
$ws2.Range("C6").Value = 65618880.5129077
$ws2.Range("D6").Value = 57920929.21961301
$ws2.Range("E6").Value = 54612924.7336
$ws2.Range("F6").Value = 73219718

# Row 7: PassivoTotal
$ws2.Range("A7").Value = "PassivoTotal"
$ws2.Range("B7").Value = 69691047.31552862
$ws2.Range("C7").Value = 65618880.5129077
$ws2.Range("D7").Value = 57920929.21961301
$ws2.Range("E7").Value = 54612924.7336
$ws2.Range("F7").Value = 73219718

# Row 8: PassivoCirculante
$ws2.Range("A8").Value = "PassivoCirculante"
$ws2.Range("B8").Value = 16084627.91634959
$ws2.Range("C8").Value = 18279637.44694333
$ws2.Range("D8").Value = 18373669.44346336
$ws2.Range("E8").Value = 18565945.7574
$ws2.Range("F8").Value = 27347016

# Row 9: PassivoNaoCirculante
$ws2.Range("A9").Value = "PassivoNaoCirculante"
$ws2.Range("B9").Value = 34305057.82322473
$ws2.Range("C9").Value = 29409556.67532437
$ws2.Range("D9").Value = 23417345.436912
$ws2.Range("E9").Value = 20652951.5502
$ws2.Range("F9").Value = 25161794

# Row 10: PatrimonioLiquido
$ws2.Range("A10").Value = "PatrimonioLiquido"
$ws2.Range("B10").Value = 19301361.5759543
$ws2.Range("C10").Value = 17929686.39064
$ws2.Range("D10").Value = 16129914.33923764
$ws2.Range("E10").Value = 15394027.426
$ws2.Range("F10").Value = 20710908

# Row 11: ReceitaLiquida
$ws2.Range("A11").Value = "ReceitaLiquida"
$ws2.Range("B11").Value = 22827248.50145473
$ws2.Range("C11").Value = 25011789.94273384
$ws2.Range("D11").Value = 24647341.8045608
$ws2.Range("E11").Value = 27419853.60827202
$ws2.Range("F11").Value = 35424174

# Row 12: LucroBruto
$ws2.Range("A12").Value = "LucroBruto"
$ws2.Range("B12").Value = 2797304.76988113
$ws2.Range("C12").Value = 3904740.967449251
$ws2.Range("D12").Value = 4950725.889129187
$ws2.Range("E12").Value = 4729220.957990627
$ws2.Range("F12").Value = 6382234

# Row 13: LucroOperacional
$ws2.Range("A13").Value = "LucroOperacional"
$ws2.Range("B13").Value = -1953707.982647893
$ws2.Range("C13").Value = 1163684.532499225
$ws2.Range("D13").Value = -583986.5696378001
$ws2.Range("E13").Value = 1598961.670603361
$ws2.Range("F13").Value = 3772121

# Row 14: LucroLiquido
$ws2.Range("A14").Value = "LucroLiquido"
$ws2.Range("B14").Value = -4179466.994926325
$ws2.Range("C14").Value = -295266.5330622761
$ws2.Range("D14").Value = -1100740.706332947
$ws2.Range("E14").Value = 823687.2244588722
$ws2.Range("F14").Value = 1923831

# Row 15: AtivoPermanente
$ws2.Range("A15").Value = "AtivoPermanente"
$ws2.Range("B15").Value = 29141300.63507662
$ws2.Range("C15").Value = 27650888.97605425
$ws2.Range("D15").Value = 25044728.79473388
$ws2.Range("E15").Value = 23942405.2206
$ws2.Range("F15").Value = 32892812

# Row 16 ("AtivoPermanente") no longer exists after the reshuffle above —
# its data now lives in row 15, so remove the stale trailing row completely
# (including formatting) so it drops out of the sheet's used range.
$ws2.Range("A16:F16").Clear()
